# Update stats for 2025-09 (row 22 of Sheet1)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B22").Value = 6286
$ws.Range("D22").Value = 5748562
$ws.Range("E22").Value = 914.5023862551702
$ws.Range("F22").Value = 8.211396109485293
$ws.Range("H22").Value = 25.01330909280699
